# Auto-generated edit script applying the diff to cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.189.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.09%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.932.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.12%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'592.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.93%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'146.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.69%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.03%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.507"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.40%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'6.91"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.15%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.01%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.442"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.15%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.31%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'33.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.00%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -0.56%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.417.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.10%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'61.117.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.09%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'6.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.34%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.917.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.70%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'432.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.47%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'13.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.25%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.685"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.99%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'7.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.23%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'81.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.84%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.07%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.99%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'12.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.96%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.05%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +6.58%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.01%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.15%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.49%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'26.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.01%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.109"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.99%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.0₃0864"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.23%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.92%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.67%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +3.99%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +1.46%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +2.50%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.44%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'8.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.58%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.290"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.36%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'39.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -4.79%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'376.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.02%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.61%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.717.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.45%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'131.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.32%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D49").Value = "'24.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.03%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.106"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.24%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -2.56%  "
$ws.Range("E51").Style = "Normal"
